$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 2
$ws.Range("C8").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("C19").Value = 2
$ws.Range("C21").Value = 5
$ws.Range("C22").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("C24").Value = 1
$ws.Range("C26").Value = 4
$ws.Range("C29").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("C31").Value = 1
$ws.Range("C32").Value = 1
$ws.Range("C33").Value = 1
$ws.Range("C34").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("C36").Value = 1
$ws.Range("C37").Value = 1
$ws.Range("C38").Value = 1
$ws.Range("C39").Value = 1
$ws.Range("C40").Value = 1
$ws.Range("C41").Value = 1
$ws.Range("C42").Value = 1
$ws.Range("C43").Value = 1
$ws.Range("C44").Value = 1
$ws.Range("C45").Value = 1
$ws.Range("C46").Value = 1
$ws.Range("C47").Value = 1
$ws.Range("C48").Value = 6
$ws.Range("C49").Value = 1
$ws.Range("C50").Value = 12
$ws.Range("C55").Value = 47
$ws.Range("C56").Value = 7
